$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.016579125216054
$ws.Range("D2").Value = 1.048426908005787
$ws.Range("E2").Value = 1.018066743734281
$ws.Range("F2").Value = 1.049588973064432
$ws.Range("I2").Value = 1.038683989779495
$ws.Range("J2").Value = 1.021798131300227
$ws.Range("K2").Value = 1.051186784381079
$ws.Range("L2").Value = 1.020914063343212
$ws.Range("M2").Value = 1.052345611504537
$ws.Range("N2").Value = 1.011289590342848

$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.018159180361971
$ws.Range("D3").Value = 1.049332845072452
$ws.Range("E3").Value = 1.019422828501661
$ws.Range("F3").Value = 1.050843941809018
$ws.Range("I3").Value = 1.038937344410625
$ws.Range("J3").Value = 1.023010383082386
$ws.Range("K3").Value = 1.051904521550316
$ws.Range("L3").Value = 1.022074361519569
$ws.Range("M3").Value = 1.053411714050287
$ws.Range("N3").Value = 1.011696015993448

$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.019178858726698
$ws.Range("D4").Value = 1.049914311756314
$ws.Range("E4").Value = 1.020298313975961
$ws.Range("F4").Value = 1.051650574226001
$ws.Range("I4").Value = 1.039097368684216
$ws.Range("J4").Value = 1.023791800922673
$ws.Range("K4").Value = 1.052363590475863
$ws.Range("L4").Value = 1.022822646097094
$ws.Range("M4").Value = 1.054095594038125
$ws.Range("N4").Value = 1.011957819778471

$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.019606890452698
$ws.Range("D5").Value = 1.050157629500328
$ws.Range("E5").Value = 1.020665899277383
$ws.Range("F5").Value = 1.051988391389344
$ws.Range("I5").Value = 1.039163706952734
$ws.Range("J5").Value = 1.024119601083609
$ws.Range("K5").Value = 1.052555303510983
$ws.Range("L5").Value = 1.023136632410743
$ws.Range("M5").Value = 1.054381675216878
$ws.Range("N5").Value = 1.012067602075899

$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.019678721582711
$ws.Range("D6").Value = 1.050198417383637
$ws.Range("E6").Value = 1.020727591210107
$ws.Range("F6").Value = 1.052045036824104
$ws.Range("I6").Value = 1.039174790579778
$ws.Range("J6").Value = 1.024174598889784
$ws.Range("K6").Value = 1.052587418008582
$ws.Range("L6").Value = 1.023189317561475
$ws.Range("M6").Value = 1.054429626255681
$ws.Range("N6").Value = 1.012086018663336

$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.019184580611537
$ws.Range("D7").Value = 1.049917567419246
$ws.Range("E7").Value = 1.020303227499515
$ws.Range("F7").Value = 1.051655093220651
$ws.Range("I7").Value = 1.039098258775609
$ws.Range("J7").Value = 1.023796183774998
$ws.Range("K7").Value = 1.052366157178847
$ws.Range("L7").Value = 1.022826843915513
$ws.Range("M7").Value = 1.054099422249298
$ws.Range("N7").Value = 1.011959287791448

$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.017113682618306
$ws.Range("D8").Value = 1.04873405583952
$ws.Range("E8").Value = 1.018525456312781
$ws.Range("F8").Value = 1.050014220514243
$ws.Range("I8").Value = 1.038770422950121
$ws.Range("J8").Value = 1.022208441409066
$ws.Range("K8").Value = 1.051430457338792
$ws.Range("L8").Value = 1.021306714902362
$ws.Range("M8").Value = 1.052707142776246
$ws.Range("N8").Value = 1.011427189902839

$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.013443150269472
$ws.Range("D9").Value = 1.046612159031647
$ws.Range("E9").Value = 1.015377171470738
$ws.Range("F9").Value = 1.047081051185391
$ws.Range("I9").Value = 1.038162720731101
$ws.Range("J9").Value = 1.019387374429421
$ws.Range("K9").Value = 1.049740502809432
$ws.Range("L9").Value = 1.018608524540014
$ws.Range("M9").Value = 1.050207899871445
$ws.Range("N9").Value = 1.010480400711073

$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.010981082883099
$ws.Range("D10").Value = 1.045172889142821
$ws.Range("E10").Value = 1.013267287629456
$ws.Range("F10").Value = 1.045097159644038
$ws.Range("I10").Value = 1.037737336935523
$ws.Range("J10").Value = 1.017490509755406
$ws.Range("K10").Value = 1.048586031051329
$ws.Range("L10").Value = 1.016796123454417
$ws.Range("M10").Value = 1.048510566938145
$ws.Range("N10").Value = 1.009842873861668

$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.009911255083926
$ws.Range("D11").Value = 1.04454376874548
$ws.Range("E11").Value = 1.012350950902627
$ws.Range("F11").Value = 1.044231276273299
$ws.Range("I11").Value = 1.037548320392061
$ws.Range("J11").Value = 1.016665198545045
$ws.Range("K11").Value = 1.048079485444751
$ws.Range("L11").Value = 1.016007998004312
$ws.Range("M11").Value = 1.047768130962407
$ws.Range("N11").Value = 1.009565276279977

$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.009513298364739
$ws.Range("D12").Value = 1.044309193374891
$ws.Range("E12").Value = 1.012010160169225
$ws.Range("F12").Value = 1.043908612074235
$ws.Range("I12").Value = 1.037477384932274
$ws.Range("J12").Value = 1.016358036593569
$ws.Range("K12").Value = 1.04789032830757
$ws.Range("L12").Value = 1.015714741252146
$ws.Range("M12").Value = 1.047491226242375
$ws.Range("N12").Value = 1.009461928739331

$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.009598687675607
$ws.Range("D13").Value = 1.044359551030349
$ws.Range("E13").Value = 1.012083280224056
$ws.Range("F13").Value = 1.043977871679991
$ws.Range("I13").Value = 1.037492633725073
$ws.Range("J13").Value = 1.016423951401251
$ws.Range("K13").Value = 1.047930948637344
$ws.Range("L13").Value = 1.015777669138151
$ws.Range("M13").Value = 1.047550674540469
$ws.Range("N13").Value = 1.009484107847394

$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.009878371667949
$ws.Range("D14").Value = 1.044524396896719
$ws.Range("E14").Value = 1.012322789713077
$ws.Range("F14").Value = 1.044204625943448
$ws.Range("I14").Value = 1.037542471671668
$ws.Range("J14").Value = 1.016639820836198
$ws.Range("K14").Value = 1.048063870153726
$ws.Range("L14").Value = 1.015983767807143
$ws.Range("M14").Value = 1.047745265045301
$ws.Range("N14").Value = 1.009556738362958

$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.010050617594887
$ws.Range("D15").Value = 1.044625845531142
$ws.Range("E15").Value = 1.012470303119842
$ws.Range("F15").Value = 1.04434419912688
$ws.Range("I15").Value = 1.037573082179268
$ws.Range("J15").Value = 1.016772744692514
$ws.Range("K15").Value = 1.048145634378839
$ws.Range("L15").Value = 1.016110683884953
$ws.Range("M15").Value = 1.047865008597108
$ws.Range("N15").Value = 1.00960145711917

$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.01105200394541
$ws.Range("D16").Value = 1.045214516950465
$ws.Range("E16").Value = 1.013328043206031
$ws.Range("F16").Value = 1.045154480560698
$ws.Range("I16").Value = 1.037749779586757
$ws.Range("J16").Value = 1.017545198718196
$ws.Range("K16").Value = 1.048619508316434
$ws.Range("L16").Value = 1.016848357487196
$ws.Range("M16").Value = 1.048559681771251
$ws.Range("N16").Value = 1.009861264271573

$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.011679137208786
$ws.Range("D17").Value = 1.045582190321524
$ws.Range("E17").Value = 1.013865339234986
$ws.Range("F17").Value = 1.045660910357131
$ws.Range("I17").Value = 1.037859324850243
$ws.Range("J17").Value = 1.018028672337141
$ws.Range("K17").Value = 1.048914973153912
$ws.Range("L17").Value = 1.017310178754816
$ws.Range("M17").Value = 1.048993424373567
$ws.Range("N17").Value = 1.010023818588139

$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.012044573861001
$ws.Range("D18").Value = 1.045796078355728
$ws.Range("E18").Value = 1.014178470957222
$ws.Range("F18").Value = 1.045955642387425
$ws.Range("I18").Value = 1.037922755568029
$ws.Range("J18").Value = 1.018310293295524
$ws.Range("K18").Value = 1.049086671295521
$ws.Range("L18").Value = 1.017579229306709
$ws.Range("M18").Value = 1.049245697922484
$ws.Range("N18").Value = 1.010118485008106

$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.012169117737373
$ws.Range("D19").Value = 1.045868912165209
$ws.Range("E19").Value = 1.014285196320637
$ws.Range("F19").Value = 1.046056026704364
$ws.Range("I19").Value = 1.037944304952437
$ws.Range("J19").Value = 1.018406254443132
$ws.Range("K19").Value = 1.049145107239059
$ws.Range("L19").Value = 1.017670914263418
$ws.Range("M19").Value = 1.049331594573954
$ws.Range("N19").Value = 1.010150738681128

$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.011611889035551
$ws.Range("D20").Value = 1.04554280140004
$ws.Range("E20").Value = 1.013807719833439
$ws.Range("F20").Value = 1.045606643548063
$ws.Range("I20").Value = 1.037847619808685
$ws.Range("J20").Value = 1.017976839678255
$ws.Range("K20").Value = 1.048883338956492
$ws.Range("L20").Value = 1.017260663075856
$ws.Range("M20").Value = 1.048946962543339
$ws.Range("N20").Value = 1.010006393460361

$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.009796027688017
$ws.Range("D21").Value = 1.044475878560149
$ws.Range("E21").Value = 1.012252271884427
$ws.Range("F21").Value = 1.044137881145274
$ws.Range("I21").Value = 1.037527815702626
$ws.Range("J21").Value = 1.016576269441508
$ws.Range("K21").Value = 1.048024755817462
$ws.Range("L21").Value = 1.015923091093089
$ws.Range("M21").Value = 1.047687994270247
$ws.Range("N21").Value = 1.009535357013486

$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.008650989654894
$ws.Range("D22").Value = 1.043799898340062
$ws.Range("E22").Value = 1.01127185048166
$ws.Range("F22").Value = 1.043208410498453
$ws.Range("I22").Value = 1.037322539081696
$ws.Range("J22").Value = 1.015692172005721
$ws.Range("K22").Value = 1.04747912211346
$ws.Range("L22").Value = 1.015079140172597
$ws.Range("M22").Value = 1.046889883930569
$ws.Range("N22").Value = 1.009237833788866

$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.009258316621644
$ws.Range("D23").Value = 1.044158739257282
$ws.Range("E23").Value = 1.011791826102602
$ws.Range("F23").Value = 1.04370171211901
$ws.Range("I23").Value = 1.037431759119423
$ws.Range("J23").Value = 1.016161184319306
$ws.Range("K23").Value = 1.047768924880145
$ws.Range("L23").Value = 1.015526818906237
$ws.Range("M23").Value = 1.047313600370997
$ws.Range("N23").Value = 1.009395686922693

$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.011642276701405
$ws.Range("D24").Value = 1.04556060131866
$ws.Range("E24").Value = 1.013833756376501
$ws.Range("F24").Value = 1.045631166422588
$ws.Range("I24").Value = 1.03785291025069
$ws.Range("J24").Value = 1.018000261803083
$ws.Range("K24").Value = 1.048897635071148
$ws.Range("L24").Value = 1.017283038074752
$ws.Range("M24").Value = 1.048967958873196
$ws.Range("N24").Value = 1.010014267585239

$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.014394671255187
$ws.Range("D25").Value = 1.047165053521213
$ws.Range("E25").Value = 1.016192983420399
$ws.Range("F25").Value = 1.047844330720566
$ws.Range("I25").Value = 1.038323389902831
$ws.Range("J25").Value = 1.020119496828369
$ws.Range("K25").Value = 1.050182288860626
$ws.Range("L25").Value = 1.019308434350678
$ws.Range("M25").Value = 1.050859483041291
$ws.Range("N25").Value = 1.010726272032547
